# ----------------------------------------------------------------------------
# Adds a new "2022-Q4" quarterly sheet (inserted right after "总计" and before
# the existing "2022-Q3" sheet) populated with fund-holding data, and updates
# the "总计" (summary) sheet with a new leading row for 2022-Q4 while shifting
# all the previously existing rows down by one.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a 2022-Q4 row at the top of
#    the data (row 2) and push every existing data row down by one.
# -----------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Give the brand-new row 7 the same look (border/font/alignment) as the
# existing index cells in column A before we start shuffling values around.
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial(-4122)

# Push the existing 5 data rows (old rows 2..6) down to rows 3..7,
# working from the bottom up so that we never overwrite data before
# it has been copied.
for ($i = 5; $i -ge 2; $i--) {
    $srcRow = $i
    $dstRow = $i + 1
    $summary.Cells.Item($dstRow, 1).Value = $summary.Cells.Item($srcRow, 1).Value
    $summary.Cells.Item($dstRow, 2).Value = $summary.Cells.Item($srcRow, 2).Value
    $summary.Cells.Item($dstRow, 3).Value = $summary.Cells.Item($srcRow, 3).Value
    $summary.Cells.Item($dstRow, 4).Value = $summary.Cells.Item($srcRow, 4).Value
}

# Fill in the brand-new 2022-Q4 summary row.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 23
$summary.Range("D2").Value = 3.01

# -----------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right before the existing
#    "2022-Q3" sheet (i.e. as the second sheet overall).
# -----------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($oldQ3)
$newSheet.Name = "2022-Q4"

# Copy the header-row formatting (bold/centered/bordered style) from the
# neighbouring sheet so the new sheet matches the look of the others.
$oldQ3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the index-column (A) formatting down for all 23 data rows.
$oldQ3.Range("A2").Copy()
$newSheet.Range("A2:A24").PasteSpecial(-4122)

# Columns B (fund code) and D:G (numeric-looking figures) are stored as
# plain text in the source data (e.g. "005344", "0.3550" keeps the
# trailing zero), so mark those ranges as text before writing the values
# to stop Excel from re-interpreting them as numbers.
$newSheet.Range("B2:B24").NumberFormat = "@"
$newSheet.Range("D2:G24").NumberFormat = "@"

$fundRows = @(
    @(0, '005344', '长安裕盛灵活配置混合C', '3.88', '91.94', '9.15', '0.3550', 2),
    @(1, '001239', '长盛国企改革主题灵活配置混合', '4.04', '85.43', '8.39', '0.3390', 3),
    @(2, '005478', '长安鑫禧灵活配置混合C', '3.75', '92.51', '9.01', '0.3379', 3),
    @(3, '501075', '万家科创主题灵活配置混合（LOF）A', '7.37', '94.25', '4.20', '0.3095', 5),
    @(4, '010694', '万家内需增长一年持有期混合', '9.46', '94.46', '3.05', '0.2885', 9),
    @(5, '005119', '银华智荟内在价值灵活配置混合A', '6.43', '94.86', '4.16', '0.2675', 8),
    @(6, '009859', '银华乐享混合A', '4.40', '94.60', '5.54', '0.2438', 7),
    @(7, '010611', '万家战略发展产业混合A', '5.75', '92.07', '3.60', '0.2070', 8),
    @(8, '010612', '万家战略发展产业混合C', '4.14', '92.07', '3.60', '0.1490', 8),
    @(9, '005477', '长安鑫禧灵活配置混合A', '1.08', '92.51', '9.01', '0.0973', 3),
    @(10, '290014', '泰信现代服务业混合', '2.40', '77.48', '3.27', '0.0785', 9),
    @(11, '015687', '银华乐享混合C', '1.29', '94.60', '5.54', '0.0715', 7),
    @(12, '013842', '银华新锐成长混合A', '1.80', '94.70', '3.61', '0.0650', 9),
    @(13, '005343', '长安裕盛灵活配置混合A', '0.54', '91.94', '9.15', '0.0494', 2),
    @(14, '290008', '泰信发展主题混合', '1.18', '84.83', '3.44', '0.0406', 8),
    @(15, '016262', '银华智荟内在价值灵活配置混合C', '0.82', '94.86', '4.16', '0.0341', 8),
    @(16, '000354', '长盛城镇化主题混合', '0.33', '85.83', '8.17', '0.0270', 2),
    @(17, '013843', '银华新锐成长混合C', '0.64', '94.70', '3.61', '0.0231', 9),
    @(18, '005186', '长安鑫兴灵活配置混合A', '0.34', '92.75', '6.12', '0.0208', 9),
    @(19, '005187', '长安鑫兴灵活配置混合C', '0.07', '92.75', '6.12', '0.0043', 9),
    @(20, '007501', '万家科创主题灵活配置混合（LOF）C', '0.10', '94.25', '4.20', '0.0042', 5),
    @(21, '011987', '财通资管智选核心回报6个月持有期混合A', '0.11', '38.44', '1.23', '0.0014', 7),
    @(22, '011988', '财通资管智选核心回报6个月持有期混合C', '0.01', '38.44', '1.23', '0.0001', 7)
)

foreach ($row in $fundRows) {
    $r = $row[0] + 2
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}
